# Move the "min_weekly_freq" parameter row from the "mobility" sheet
# (where it was the last row, 26) to the "time" sheet (as new row 6),
# and update the view state so "time" becomes the active tab.

$wb = $excel.ActiveWorkbook

$mobility = $wb.Worksheets.Item("mobility")
$time = $wb.Worksheets.Item("time")

# Move A26:C26 (min_weekly_freq / 2 / description) from "mobility" to
# "time" as the new row 6, then remove the now-empty row 26 so the
# sheet's dimension shrinks back down (A1:C26 -> A1:C25).
$mobility.Range("A26:C26").Cut($time.Range("A6:C6"))
$mobility.Rows.Item(26).Delete()

# Keep the numeric cell's centered look consistent with the rest of the
# column (matches the formatting used for the equivalent cell on the
# "mobility" sheet before the move).
$time.Range("B6").HorizontalAlignment = -4108

# "time" becomes the active sheet/tab (previously "mobility" was).
$time.Activate()

# Update the on-screen selections to match where the edit left the
# cursor on each sheet.
$mobility.Range("A26:C26").Select()
$time.Range("A4").Select()

# Scroll "mobility" down a bit and "time" back to the top, mirroring the
# final view positions after the edit.
$excel.ActiveWindow.ScrollRow = 7
$mobility.Activate()
$excel.ActiveWindow.ScrollRow = 7

$time.Activate()
$time.Range("A4").Select()
